$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3"   = -7.331999999999999
    "B12"  = 5.325
    "D14"  = -7.318000000000001
    "D26"  = -8.119
    "B27"  = 6.229000000000001
    "D31"  = -8.414000000000001
    "B32"  = 6.614
    "D35"  = -7.672
    "B36"  = 8.548
    "D37"  = -7.741
    "B38"  = 5.825
    "D45"  = -7.498
    "B46"  = 6.572
    "D52"  = -7.961000000000001
    "B54"  = 5.05
    "B55"  = 4.684
    "B56"  = 4.627
    "D57"  = -8.101000000000001
    "B67"  = 4.915999999999999
    "B69"  = 5.063000000000001
    "B72"  = 5.097
    "D81"  = -6.817
    "B83"  = 5.737
    "D83"  = -8.572000000000001
    "B86"  = 5.145999999999999
    "B91"  = 6.109
    "B93"  = 5.33
    "B99"  = 5.217000000000001
    "D100" = -8.313000000000001
    "D102" = -7.861000000000002
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
